$p = $ppt.ActivePresentation

# --- Update the auto-date fields cached on the Handout Master and the
#     Notes Master (shown on printouts) from 5/27/2018 to 6/1/2018. ---
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "6/1/2018"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "6/1/2018"

# --- Slide 2: the "Lecture 6." caption becomes "Lecture 5." now that the
#     lecture has been expanded to cover the Stanford heart transplant
#     data and renumbered. ---
$s2 = $p.Slides.Item(2)
$lectureShape = $s2.Shapes.Item("Text Placeholder 5")
$lectureShape.TextFrame.TextRange.Text = "Lecture 5."
